$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Copy the formatting (number format / font / border / alignment) that
#    already lives in column I onto the five new year columns J:N. Column I
#    carries the correct style for every row that needs one (rows 3-4, the
#    data rows 6-19, and the blank filler rows 20-29). Row 5 is a section
#    header with no D:I values at all, so it is intentionally skipped — it
#    must not gain any J:N cells either.
# ---------------------------------------------------------------------------
$ws.Range("I3:I4").Copy()
$ws.Range("J3:N3").PasteSpecial(-4122)

$ws.Range("I6:I19").Copy()
$ws.Range("J6:N6").PasteSpecial(-4122)

$ws.Range("I20:I29").Copy()
$ws.Range("J20:N20").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2. Header row: years 2019-2023 in J4:N4.
# ---------------------------------------------------------------------------
$ws.Range("J4").Value2 = 2019
$ws.Range("K4").Value2 = 2020
$ws.Range("L4").Value2 = 2021
$ws.Range("M4").Value2 = 2022
$ws.Range("N4").Value2 = 2023

# ---------------------------------------------------------------------------
# 3. Data values for the new year columns, row by row.
# ---------------------------------------------------------------------------
$ws.Range("J6").Value2 = 81.1
$ws.Range("K6").Value2 = 85.8
$ws.Range("L6").Value2 = 78.1
$ws.Range("M6").Value2 = 72.2
$ws.Range("N6").Value2 = 75.7

$ws.Range("J7").Value2 = 18.9
$ws.Range("K7").Value2 = 14.2
$ws.Range("L7").Value2 = 21.9
$ws.Range("M7").Value2 = 27.8
$ws.Range("N7").Value2 = 24.3

$ws.Range("J9").Value2 = 22.8
$ws.Range("K9").Value2 = 25.6
$ws.Range("L9").Value2 = 24.2
$ws.Range("M9").Value2 = 21.4
$ws.Range("N9").Value2 = 31.1

$ws.Range("J10").Value2 = 77.2
$ws.Range("K10").Value2 = 74.4
$ws.Range("L10").Value2 = 75.8
$ws.Range("M10").Value2 = 78.6
$ws.Range("N10").Value2 = 68.9

$ws.Range("J12").Value2 = 84.4
$ws.Range("K12").Value2 = 72.7
$ws.Range("L12").Value2 = 73.3
$ws.Range("M12").Value2 = 72.8
$ws.Range("N12").Value2 = 76.7

$ws.Range("J13").Value2 = 15.6
$ws.Range("K13").Value2 = 27.3
$ws.Range("L13").Value2 = 26.7
$ws.Range("M13").Value2 = 27.2
$ws.Range("N13").Value2 = 23.3

$ws.Range("J15").Value2 = 90.3
$ws.Range("K15").Value2 = 93.4
$ws.Range("L15").Value2 = 90.5
$ws.Range("M15").Value2 = 87.8
$ws.Range("N15").Value2 = 89

$ws.Range("J16").Value2 = 9.7
$ws.Range("K16").Value2 = 6.6
$ws.Range("L16").Value2 = 9.5
$ws.Range("M16").Value2 = 12.2
$ws.Range("N16").Value2 = 11

$ws.Range("J18").Value2 = 60.2
$ws.Range("K18").Value2 = 66
$ws.Range("L18").Value2 = 59.3
$ws.Range("M18").Value2 = 44.9
$ws.Range("N18").Value2 = 48.3

$ws.Range("J19").Value2 = 39.8
$ws.Range("K19").Value2 = 34
$ws.Range("L19").Value2 = 40.7
$ws.Range("M19").Value2 = 55.1
$ws.Range("N19").Value2 = 51.7

# ---------------------------------------------------------------------------
# 4. Row 20 (the footnote row "*according to the MOF of the KR") switches to
#    the smaller 8pt Times New Roman font and gets an explicit row height.
# ---------------------------------------------------------------------------
$ws.Range("A20:C20").Font.Size = 8
$ws.Rows.Item(20).RowHeight = 15.75
